$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 15 (pushes the old rows 15-34 down to 17-36)
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(15).Insert()

# --- New row 15 ---
$ws.Cells.Item(15, 1).Value = 10
$ws.Cells.Item(15, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(15, 3).Value = "La Araucanía"
$ws.Cells.Item(15, 4).Value = 44679
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100108
$ws.Cells.Item(15, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(15, 9).Value = 100108003
$ws.Cells.Item(15, 10).Value = "Maracuyá"
$ws.Cells.Item(15, 11).Value = "Sin especificar"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 35
$ws.Cells.Item(15, 14).Value = 34000
$ws.Cells.Item(15, 15).Value = 34000
$ws.Cells.Item(15, 16).Value = 34000
$ws.Cells.Item(15, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(15, 18).Value = "Perú"
$ws.Cells.Item(15, 19).Value = 1889
$ws.Cells.Item(15, 20).Value = 18

# --- New row 16 ---
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(16, 3).Value = "La Araucanía"
$ws.Cells.Item(16, 4).Value = 44679
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100108
$ws.Cells.Item(16, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(16, 9).Value = 100108003
$ws.Cells.Item(16, 10).Value = "Maracuyá"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 55
$ws.Cells.Item(16, 14).Value = 28000
$ws.Cells.Item(16, 15).Value = 28000
$ws.Cells.Item(16, 16).Value = 28000
$ws.Cells.Item(16, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(16, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(16, 19).Value = 1556
$ws.Cells.Item(16, 20).Value = 18
